$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the end time of the last existing entry (row 85):
#    0.75 (18:00) -> 0.8125 (19:30). Dependent formulas (F85, G85) and the
#    summary rows below recalculate automatically.
$ws.Range("E85").Value = 0.8125

# 2. Insert a new row at position 86, shifting the blank separator row and
#    the three summary rows (old rows 86-89) down to rows 87-90.
$ws.Rows("86:86").Insert()

# 3. Fill in the newly inserted row 86 with the new timesheet entry
#    (2014-03-21, 20:00 -> 22:00).
$ws.Range("A86").Value = 2014
$ws.Range("B86").Value = 3
$ws.Range("C86").Value = 21
$ws.Range("D86").Value = 0.83333333333333337
$ws.Range("E86").Value = 0.91666666666666663
$ws.Range("F86").Formula = "=(E86-D86)*24*60"
$ws.Range("G86").Formula = "=F86/60"

# 4. Match the saved selection state (active cell moved down with the rows).
$ws.Range("A87").Select()
